# fix to data import script + premature rerun of analysis (before pulling)
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Country" sheet: rerun of country-level aggregation produced new
# percentages for the single data row.
# ------------------------------------------------------------------
$wsCountry = $wb.Worksheets.Item("Country")
$wsCountry.Range("A2").Value = 0.009602194787379973
$wsCountry.Range("B2").Value = 0.2112482853223594
$wsCountry.Range("C2").Value = 0.2400548696844993
$wsCountry.Range("D2").Value = 0.2235939643347051
$wsCountry.Range("E2").Value = 0.1920438957475994

# ------------------------------------------------------------------
# "States" sheet: data-import fix re-derives the per-state rows (new
# sort order + refreshed percentages) and adds an ISO STATE.CODE
# column.
# ------------------------------------------------------------------
$wsStates = $wb.Worksheets.Item("States")

# New header cell for the added column (match the existing bold,
# centered header formatting used by A1:F1).
$wsStates.Range("G1").Value = "STATE.CODE"
$wsStates.Range("G1").Font.Bold = $true
$wsStates.Range("G1").HorizontalAlignment = -4108

$statesData = @(
    @("Bihar", 0, 0.02631578947368421, 0.07894736842105263, 0.1842105263157895, 0.7631578947368421, "IN-BR"),
    @("Mizoram", 0, 0.1818181818181818, 0.3636363636363636, 0.3636363636363636, 0.5454545454545454, "IN-MZ"),
    @("Jharkhand", 0, 0.08333333333333333, 0.125, 0.3333333333333333, 0.5416666666666666, "IN-JH"),
    @("Daman and Diu", 0, 0, 0.5, 0.5, 0.5, "IN-DD"),
    @("Puducherry", 0, 0, 0, 0.25, 0.5, "IN-PY"),
    @("Uttar Pradesh", 0, 0.08, 0.1866666666666667, 0.32, 0.4533333333333333, "IN-UP"),
    @("Arunachal Pradesh", 0, 0.12, 0.28, 0.12, 0.4, "IN-AR"),
    @("Manipur", 0, 0.125, 0.25, 0.4375, 0.375, "IN-MN"),
    @("Nagaland", 0, 0.09090909090909091, 0.5454545454545454, 0.2727272727272727, 0.3636363636363636, "IN-NL"),
    @("Punjab", 0, 0.09090909090909091, 0.2272727272727273, 0.4090909090909091, 0.2727272727272727, "IN-PB"),
    @("Meghalaya", 0.09090909090909091, 0.3636363636363636, 0.3636363636363636, 0.5454545454545454, 0.1818181818181818, "IN-ML"),
    @("Odisha", 0, 0.2, 0.2666666666666667, 0.4, 0.1666666666666667, "IN-OR"),
    @("Telangana", 0, 0.2424242424242424, 0.3636363636363636, 0.2121212121212121, 0.1515151515151515, "IN-TS"),
    @("Tripura", 0, 0, 0.25, 0.625, 0.125, "IN-TR"),
    @("Assam", 0.0303030303030303, 0.2121212121212121, 0.303030303030303, 0.2424242424242424, 0.1212121212121212, "IN-AS"),
    @("Delhi", 0, 0.1818181818181818, 0.3636363636363636, 0.4545454545454545, 0.09090909090909091, "IN-DL"),
    @("Jammu and Kashmir", 0, 0.09090909090909091, 0.4545454545454545, 0.4090909090909091, 0.09090909090909091, "IN-JK"),
    @("Haryana", 0, 0.1363636363636364, 0.4090909090909091, 0.3636363636363636, 0.09090909090909091, "IN-HR"),
    @("Madhya Pradesh", 0, 0.1538461538461539, 0.4038461538461539, 0.3461538461538461, 0.07692307692307693, "IN-MP"),
    @("West Bengal", 0, 0.391304347826087, 0.3043478260869565, 0.08695652173913043, 0.04347826086956522, "IN-WB"),
    @("Chhattisgarh", 0, 0.2592592592592592, 0.4444444444444444, 0.1111111111111111, 0.03703703703703703, "IN-CT"),
    @("Rajasthan", 0.0303030303030303, 0.303030303030303, 0.303030303030303, 0.1818181818181818, 0.0303030303030303, "IN-RJ"),
    @("Himachal Pradesh", 0, 0.4166666666666667, 0.1666666666666667, 0.1666666666666667, 0, "IN-HP"),
    @("Gujarat", 0, 0.5454545454545454, 0.09090909090909091, 0.06060606060606061, 0, "IN-GJ"),
    @("Maharashtra", 0, 0.3055555555555556, 0.1944444444444444, 0.05555555555555555, 0, "IN-MH"),
    @("Karnataka", 0, 0.3, 0.06666666666666667, 0.03333333333333333, 0, "IN-KA"),
    @("Tamil Nadu", 0.05405405405405406, 0.4324324324324325, 0.05405405405405406, 0.02702702702702703, 0, "IN-TN"),
    @("Chandigarh", 0, 0, 1, 0, 0, "IN-CH"),
    @("Dadra and Nagar Haveli", 0, 0, 1, 0, 0, "IN-DN"),
    @("Andhra Pradesh", 0.07692307692307693, 0.5384615384615384, 0.1538461538461539, 0, 0, "IN-AP"),
    @("Uttarakhand", 0.07692307692307693, 0.2307692307692308, 0.07692307692307693, 0, 0, "IN-UL"),
    @("Ladakh", 0, 0.5, 0, 0, 0, "IN-LA")
)

for ($i = 0; $i -lt $statesData.Length; $i++) {
    $rowVals = $statesData[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $wsStates.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

# ------------------------------------------------------------------
# "Dark clusters" sheet: same rerun refreshed the two cluster rows.
# ------------------------------------------------------------------
$wsClusters = $wb.Worksheets.Item("Dark clusters")
$wsClusters.Range("E2").Value = 0.1777777777777778

$wsClusters.Range("B3").Value = 0.01739130434782609
$wsClusters.Range("C3").Value = 0.1652173913043478
$wsClusters.Range("D3").Value = 0.3217391304347826
$wsClusters.Range("E3").Value = 0.3130434782608696
$wsClusters.Range("F3").Value = 0.2869565217391304

Write-Host "edit applied"
